# Swap the data (columns B through AD) between specific pairs of rows.
# Column A (the running index) stays put on each row; all other fields
# (match id, date, teams, scores, odds, etc.) swap between the two rows
# in each pair, effectively reordering how the two fixtures are listed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRows  = @(23, 35, 58, 97, 131, 156, 158, 163, 166, 170, 174)
$secondRows = @(24, 36, 60, 98, 132, 157, 159, 164, 167, 171, 175)

for ($i = 0; $i -lt $firstRows.Count; $i++) {
    $r1 = $firstRows[$i]
    $r2 = $secondRows[$i]

    $range1 = $ws.Range("B$($r1):AD$($r1)")
    $range2 = $ws.Range("B$($r2):AD$($r2)")

    $values1 = $range1.Value2
    $values2 = $range2.Value2

    $range1.Value2 = $values2
    $range2.Value2 = $values1
}
